$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(15, 16),
    @(20, 21),
    @(31, 32),
    @(64, 65),
    @(74, 75),
    @(81, 82),
    @(89, 91),
    @(98, 99)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 2; $col -le 28; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Append the new match as row 100 (copy formatting from the row above first)
$newRow = 100
$ws.Cells.Item(99, 1).Copy($ws.Cells.Item($newRow, 1))
$ws.Cells.Item(99, 4).Copy($ws.Cells.Item($newRow, 4))

$ws.Cells.Item($newRow, 1).Value2 = 98
$ws.Cells.Item($newRow, 2).Value2 = 7803364
$ws.Cells.Item($newRow, 3).Value2 = "Canada Premier League"
$ws.Cells.Item($newRow, 4).Value2 = 45410.75
$ws.Cells.Item($newRow, 5).Value2 = "Cavalry FC"
$ws.Cells.Item($newRow, 6).Value2 = "HFX Wanderers"
$ws.Cells.Item($newRow, 7).Value2 = 0
$ws.Cells.Item($newRow, 8).Value2 = 0
$ws.Cells.Item($newRow, 9).Value2 = "D"
$ws.Cells.Item($newRow, 10).Value2 = 2
$ws.Cells.Item($newRow, 11).Value2 = 3.2
$ws.Cells.Item($newRow, 12).Value2 = 3.3
$ws.Cells.Item($newRow, 13).Value2 = 1.8
$ws.Cells.Item($newRow, 14).Value2 = 3.1
$ws.Cells.Item($newRow, 15).Value2 = 4.2
$ws.Cells.Item($newRow, 16).Value2 = -0.5
$ws.Cells.Item($newRow, 17).Value2 = 1.825
$ws.Cells.Item($newRow, 18).Value2 = 1.975
$ws.Cells.Item($newRow, 19).Value2 = 2.25
$ws.Cells.Item($newRow, 20).Value2 = 1.95
$ws.Cells.Item($newRow, 21).Value2 = 1.85
$ws.Cells.Item($newRow, 22).Value2 = -1
$ws.Cells.Item($newRow, 23).Value2 = 2.1
$ws.Cells.Item($newRow, 24).Value2 = -1
$ws.Cells.Item($newRow, 25).Value2 = -1
$ws.Cells.Item($newRow, 26).Value2 = 0.9750000000000001
$ws.Cells.Item($newRow, 27).Value2 = -1
$ws.Cells.Item($newRow, 28).Value2 = 0.8500000000000001

Write-Host "Canada Premier League base update applied"
